$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5218626666666667
$ws.Range("H2").Value = 1.565588
$ws.Range("I2").Value = 0.004473448591865398
$ws.Range("J2").Value = 0.004473448591865398
$ws.Range("M2").Value = 14.129345
$ws.Range("N2").Value = 42.388035
$ws.Range("O2").Value = 0.3414817166893976
$ws.Range("P2").Value = 0.3414817166893976
$ws.Range("Q2").Value = 7.373577659953334
$ws.Range("R2").Value = 66.36219893958
$ws.Range("S2").Value = 0.001527600904671965
$ws.Range("T2").Value = 0.001527600904671965

$ws.Range("G3").Value = 0.5218626666666667
$ws.Range("H3").Value = 1.565588
$ws.Range("I3").Value = 0.004473448591865398
$ws.Range("J3").Value = 0.004473448591865398
$ws.Range("O3").Value = 0.3553528814026711
$ws.Range("P3").Value = 0.3553528814026711
$ws.Range("Q3").Value = 7.673096214676889
$ws.Range("R3").Value = 69.057865932092
$ws.Range("S3").Value = 0.001589652846926091
$ws.Range("T3").Value = 0.001589652846926091

$ws.Range("G4").Value = 0.5218626666666667
$ws.Range("H4").Value = 1.565588
$ws.Range("I4").Value = 0.004473448591865398
$ws.Range("J4").Value = 0.004473448591865398
$ws.Range("O4").Value = 0.3031654019079313
$ws.Range("P4").Value = 0.3031654019079312
$ws.Range("Q4").Value = 6.546217631945335
$ws.Range("R4").Value = 58.915958687508
$ws.Range("S4").Value = 0.001356194840267343
$ws.Range("T4").Value = 0.001356194840267342

$ws.Range("I5").Value = 0.97905015906109
$ws.Range("J5").Value = 0.97905015906109
$ws.Range("M5").Value = 14.129345
$ws.Range("N5").Value = 42.388035
$ws.Range("O5").Value = 0.3414817166893976
$ws.Range("P5").Value = 0.3414817166893976
$ws.Range("Q5").Value = 1613.766702036983
$ws.Range("R5").Value = 14523.90031833285
$ws.Range("S5").Value = 0.3343277290412088
$ws.Range("T5").Value = 0.3343277290412088

$ws.Range("I6").Value = 0.97905015906109
$ws.Range("J6").Value = 0.97905015906109
$ws.Range("O6").Value = 0.3553528814026711
$ws.Range("P6").Value = 0.3553528814026711
$ws.Range("S6").Value = 0.3479082950601018
$ws.Range("T6").Value = 0.3479082950601018

$ws.Range("I7").Value = 0.97905015906109
$ws.Range("J7").Value = 0.97905015906109
$ws.Range("O7").Value = 0.3031654019079313
$ws.Range("P7").Value = 0.3031654019079312
$ws.Range("S7").Value = 0.2968141349597794
$ws.Range("T7").Value = 0.2968141349597794

$ws.Range("I8").Value = 0.0164763923470446
$ws.Range("J8").Value = 0.0164763923470446
$ws.Range("M8").Value = 14.129345
$ws.Range("N8").Value = 42.388035
$ws.Range("O8").Value = 0.3414817166893976
$ws.Range("P8").Value = 0.3414817166893976
$ws.Range("Q8").Value = 27.15800931471833
$ws.Range("R8").Value = 244.422083832465
$ws.Range("S8").Value = 0.005626386743516842
$ws.Range("T8").Value = 0.005626386743516842

$ws.Range("I9").Value = 0.0164763923470446
$ws.Range("J9").Value = 0.0164763923470446
$ws.Range("O9").Value = 0.3553528814026711
$ws.Range("P9").Value = 0.3553528814026711
$ws.Range("S9").Value = 0.005854933495643217
$ws.Range("T9").Value = 0.005854933495643217

$ws.Range("I10").Value = 0.0164763923470446
$ws.Range("J10").Value = 0.0164763923470446
$ws.Range("O10").Value = 0.3031654019079313
$ws.Range("P10").Value = 0.3031654019079312
$ws.Range("S10").Value = 0.004995072107884539
$ws.Range("T10").Value = 0.004995072107884538
